$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 298.3111
$ws.Range("I33").Value = 279.38095
$ws.Range("K33").Value = 279.38095
$ws.Range("M33").Value = -50.38094999999998

$ws.Range("H116").Value = 2229.0667
$ws.Range("J116").Value = 2648.4443
$ws.Range("L116").Value = 2648.4443
$ws.Range("N116").Value = -9532.444299999999

$ws.Range("H129").Value = 919.1579
$ws.Range("I129").Value = 423.7
$ws.Range("J129").Value = 994.2273
$ws.Range("K129").Value = 1271.1
$ws.Range("L129").Value = 2982.6819
$ws.Range("M129").Value = 3728.9
$ws.Range("N129").Value = -12982.6819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2314.923
$ws.Range("I2").Value = 1231.2222
$ws.Range("J2").Value = 4753.25
$ws.Range("K2").Value = 1231.2222
$ws.Range("L2").Value = 4753.25
$ws.Range("M2").Value = -1118.2222
$ws.Range("N2").Value = -4979.25

$ws.Range("H32").Value = 17638.09
$ws.Range("I32").Value = 19064.45
$ws.Range("K32").Value = 19064.45
$ws.Range("M32").Value = -18777.45

$ws.Range("H45").Value = 3281.2273
$ws.Range("I45").Value = 2390.375
$ws.Range("J45").Value = 3790.2856
$ws.Range("K45").Value = 2390.375
$ws.Range("L45").Value = 3790.2856
$ws.Range("M45").Value = -2013.375
$ws.Range("N45").Value = -4544.2856

$ws.Range("H63").Value = 2842418
$ws.Range("I63").Value = 1777.5555
$ws.Range("J63").Value = 15625300
$ws.Range("K63").Value = 1777.5555
$ws.Range("L63").Value = 15625300
$ws.Range("M63").Value = -1091.5555
$ws.Range("N63").Value = -15626672

$ws.Range("H66").Value = 2842418
$ws.Range("I66").Value = 1777.5555
$ws.Range("J66").Value = 15625300
$ws.Range("K66").Value = 8887.7775
$ws.Range("L66").Value = 78126500
$ws.Range("M66").Value = -5455.7775
$ws.Range("N66").Value = -78133364

$ws.Range("H88").Value = 52774.6
$ws.Range("I88").Value = 1411
$ws.Range("J88").Value = 74787.57000000001
$ws.Range("K88").Value = 1411
$ws.Range("L88").Value = 74787.57000000001
$ws.Range("M88").Value = -1005
$ws.Range("N88").Value = -75599.57000000001

$ws.Range("H91").Value = 52774.6
$ws.Range("I91").Value = 1411
$ws.Range("J91").Value = 74787.57000000001
$ws.Range("K91").Value = 1411
$ws.Range("L91").Value = 74787.57000000001
$ws.Range("M91").Value = -7
$ws.Range("N91").Value = -77595.57000000001

$ws.Range("H116").Value = 2314.923
$ws.Range("I116").Value = 1231.2222
$ws.Range("J116").Value = 4753.25
$ws.Range("K116").Value = 1231.2222
$ws.Range("L116").Value = 4753.25
$ws.Range("M116").Value = 1062.7778
$ws.Range("N116").Value = -9341.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2314.923
$ws.Range("I3").Value = 1231.2222
$ws.Range("J3").Value = 4753.25
$ws.Range("K3").Value = 1231.2222
$ws.Range("L3").Value = 4753.25
$ws.Range("M3").Value = -1117.2222
$ws.Range("N3").Value = -4981.25

$ws.Range("H99").Value = 2333.3333
$ws.Range("I99").Value = 2333.3333
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2333.3333
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -835.3332999999998
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 1284.375
$ws.Range("I107").Value = 1303.3334
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1303.3334
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 616.6666
$ws.Range("N107").Value = -4840

$ws.Range("H134").Value = 56881.74
$ws.Range("I134").Value = 63277.94
$ws.Range("J134").Value = 2514
$ws.Range("K134").Value = 189833.82
$ws.Range("L134").Value = 7542
$ws.Range("M134").Value = -187298.82
$ws.Range("N134").Value = -12612

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10882.724
$ws.Range("I31").Value = 19137.908
$ws.Range("J31").Value = 3618.16
$ws.Range("K31").Value = 19137.908
$ws.Range("L31").Value = 3618.16
$ws.Range("M31").Value = -18842.908
$ws.Range("N31").Value = -4208.16

$ws.Range("H34").Value = 10882.724
$ws.Range("I34").Value = 19137.908
$ws.Range("J34").Value = 3618.16
$ws.Range("K34").Value = 19137.908
$ws.Range("L34").Value = 3618.16
$ws.Range("M34").Value = -18935.908
$ws.Range("N34").Value = -4022.16

$ws.Range("H58").Value = 14607.838
$ws.Range("I58").Value = 1148.0646
$ws.Range("J58").Value = 84150
$ws.Range("K58").Value = 1148.0646
$ws.Range("L58").Value = 84150
$ws.Range("M58").Value = -945.0645999999999
$ws.Range("N58").Value = -84556

$ws.Range("H94").Value = 3645.353
$ws.Range("I94").Value = 2037.125
$ws.Range("J94").Value = 5074.8887
$ws.Range("K94").Value = 2037.125
$ws.Range("L94").Value = 5074.8887
$ws.Range("M94").Value = -1586.125
$ws.Range("N94").Value = -5976.8887

$ws.Range("H132").Value = 20715.139
$ws.Range("I132").Value = 27062.45
$ws.Range("K132").Value = 81187.35000000001
$ws.Range("M132").Value = -78657.35000000001

$ws.Range("H136").Value = 14607.838
$ws.Range("I136").Value = 1148.0646
$ws.Range("J136").Value = 84150
$ws.Range("K136").Value = 3444.1938
$ws.Range("L136").Value = 252450
$ws.Range("M136").Value = -894.1938
$ws.Range("N136").Value = -257550

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2880
$ws.Range("J80").Value = 2880
$ws.Range("L80").Value = 8640
$ws.Range("N80").Value = -10512

$ws.Range("H83").Value = 2880
$ws.Range("J83").Value = 2880
$ws.Range("L83").Value = 25920
$ws.Range("N83").Value = -35280

$ws.Range("H112").Value = 2475
$ws.Range("I112").Value = 950
$ws.Range("J112").Value = 4000
$ws.Range("K112").Value = 2850
$ws.Range("L112").Value = 12000
$ws.Range("M112").Value = -1742
$ws.Range("N112").Value = -14216

$ws.Range("H131").Value = 110697.51
$ws.Range("J131").Value = 121299.56
$ws.Range("L131").Value = 363898.68
$ws.Range("N131").Value = -373978.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 11262.75
$ws.Range("J92").Value = 11262.75
$ws.Range("L92").Value = 11262.75
$ws.Range("N92").Value = -15006.75

$ws.Range("H122").Value = 3165.5715
$ws.Range("I122").Value = 3192.5
$ws.Range("J122").Value = 3004
$ws.Range("K122").Value = 9577.5
$ws.Range("L122").Value = 9012
$ws.Range("M122").Value = -7127.5
$ws.Range("N122").Value = -13912

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 274
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H82").Value = 3583.3333
$ws.Range("I82").Value = 3900
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 3900
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -3539
$ws.Range("N82").Value = -2722

$ws.Range("H85").Value = 3583.3333
$ws.Range("I85").Value = 3900
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 3900
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -2652
$ws.Range("N85").Value = -4496

$ws.Range("H93").Value = 3624.125
$ws.Range("I93").Value = 3665.5
$ws.Range("K93").Value = 3665.5
$ws.Range("M93").Value = -2417.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H125").Value = 53936.668
$ws.Range("J125").Value = 53936.668
$ws.Range("L125").Value = 53936.668
$ws.Range("N125").Value = -63776.668
